# Optimized algorithm for scaleUp & scaleDown
# This change reorders several data rows on the sheet (moving whole rows
# of A:F values to different row positions), while leaving other rows
# (e.g. header, row 2, 4, 10, 16, 18, 24, 25, 26, ...) untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping: destination row -> source row (values currently sitting in
# source row get copied into destination row). Snapshot all source rows
# first so that overlapping writes (cycles) don't clobber data that is
# still needed.
$mapping = @{
    3  = 11
    5  = 12
    6  = 13
    7  = 3
    8  = 7
    9  = 14
    11 = 15
    12 = 5
    13 = 9
    14 = 8
    15 = 6
    17 = 20
    19 = 21
    20 = 17
    21 = 19
    22 = 23
    23 = 22
}

$cols = @(1, 2, 3, 4, 5, 6)  # columns A..F

# Snapshot current values of every row that is used as a source.
$snapshot = @{}
foreach ($destRow in $mapping.Keys) {
    $srcRow = $mapping[$destRow]
    if (-not $snapshot.ContainsKey($srcRow)) {
        $rowVals = @{}
        foreach ($col in $cols) {
            $rowVals[$col] = $ws.Cells.Item($srcRow, $col).Value2
        }
        $snapshot[$srcRow] = $rowVals
    }
}

# Now write snapshot values into the destination rows.
foreach ($destRow in $mapping.Keys) {
    $srcRow = $mapping[$destRow]
    $rowVals = $snapshot[$srcRow]
    foreach ($col in $cols) {
        $ws.Cells.Item($destRow, $col).Value = $rowVals[$col]
    }
}
